$d = $word.ActiveDocument

# Mapping of old text -> new text for this revision of the worksheet.
$replacements = [ordered]@{
    "2025-07-05 Saturday" = "2025-07-06 Sunday"
    "524÷7="              = "947÷2="
    "923÷7="              = "768÷8="
    "333÷5="              = "331÷7="
    "428÷4="              = "927÷7="
    "398÷6="              = "352÷5="
    "152÷5="              = "354÷5="
    "766÷8="              = "557÷8="
    "856÷6="              = "581÷5="
    "884÷7="              = "350÷4="
    "238÷8="              = "614÷3="
    "588÷6="              = "589÷9="
    "474÷5="              = "162÷6="
    "153÷6="              = "104÷9="
    "842÷8="              = "632÷7="
    "680÷8="              = "546÷5="
    "564÷9="              = "760÷6="
    "951÷7="              = "647÷8="
    "778÷8="              = "855÷8="
    "846÷5="              = "369÷7="
    "623÷8="              = "137÷6="
    "861÷9="              = "781÷2="
    "466÷3="              = "970÷2="
    "413÷9="              = "344÷7="
    "345÷7="              = "379÷7="
    "736÷9="              = "800÷9="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
